# Rewrites the "Insights" table (A1:D37) to the updated agility-topic dataset:
#  - trims the Artea GDP-forecast insight text (D3)
#  - drops the two superseded Citadele financial-results insights and replaces them
#    with the refinancing-surge insight, shifting the remaining Citadele/Luminor/SEB rows up
#  - drops the two superseded SEB fraud/maintenance insights and appends three new SEB rows
#    (GDP forecast, Latvia expansion, Baltic merger) at the end of the table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Insights")

$ws.Cells.Item(1, 1).Value = 'Brand'
$ws.Cells.Item(1, 2).Value = 'Topic'
$ws.Cells.Item(1, 3).Value = 'Example #'
$ws.Cells.Item(1, 4).Value = 'One-line Insight'

$ws.Cells.Item(2, 1).Value = 'Artea'
$ws.Cells.Item(2, 2).Value = 'Economic Growth and Trends'
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 'Šiaulių bank’s chief economist comments on U.S. tariff policy, stressing risks from unpredictable presidential decisions.'

$ws.Cells.Item(3, 1).Value = 'Artea'
$ws.Cells.Item(3, 2).Value = 'Economic Growth and Trends'
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 'Artea forecasts GDP growth of 2.7% in 2025, slowing further to 2.5% in 2026.'

$ws.Cells.Item(4, 1).Value = 'Artea'
$ws.Cells.Item(4, 2).Value = 'Brand Development and Corporate Identity'
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 'Šiaulių bankas will rebrand as Artea in May following unanimous shareholder approval.'

$ws.Cells.Item(5, 1).Value = 'Artea'
$ws.Cells.Item(5, 2).Value = 'Brand Development and Corporate Identity'
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 'Šiaulių bankas announces a name change to Artea on May 5, with temporary service disruptions during system updates.'

$ws.Cells.Item(6, 1).Value = 'Artea'
$ws.Cells.Item(6, 2).Value = 'Banking Sector Developments and Innovations'
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 'Funds transfer from Poland’s PKO to Artea stalled when Artea requested missing documentation, leaving €45,000 unresolved.'

$ws.Cells.Item(7, 1).Value = 'Artea'
$ws.Cells.Item(7, 2).Value = 'Banking Sector Developments and Innovations'
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 'Šiaulių bankas Group posted €17.7m net profit in Q1, down 21% year-on-year.'

$ws.Cells.Item(8, 1).Value = 'Artea'
$ws.Cells.Item(8, 2).Value = 'Banking Sector Developments and Innovations'
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 'Šiaulių bankas plans to repurchase up to 7m shares, cancel 10.6m previously acquired, and may allocate some to staff.'

$ws.Cells.Item(9, 1).Value = 'Citadele'
$ws.Cells.Item(9, 2).Value = 'Economic Trends and Conditions'
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 'From June 1, Citadele will raise fees for private clients, increasing monthly card charges from €1.99 to €2.49.'

$ws.Cells.Item(10, 1).Value = 'Citadele'
$ws.Cells.Item(10, 2).Value = 'Economic Trends and Conditions'
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 'Citadele plans to adjust service fees for clients, though most banks in Lithuania are not following suit.'

$ws.Cells.Item(11, 1).Value = 'Citadele'
$ws.Cells.Item(11, 2).Value = 'Economic Trends and Conditions'
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 'Citadele reports Lithuanian firms are lowering expectations, while Estonian business sentiment hits a 29-month high.'

$ws.Cells.Item(12, 1).Value = 'Citadele'
$ws.Cells.Item(12, 2).Value = 'Financial Literacy and Management'
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 'Citadele’s surge in refinancing—from 7% to 32% of mortgage applications—shows that new rules are prompting consumers to review loans more actively, strengthening financial literacy but straining loan processing.'

$ws.Cells.Item(13, 1).Value = 'Citadele'
$ws.Cells.Item(13, 2).Value = 'Financial Literacy and Management'
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = 'A Citadele survey shows Lithuanians most often save for specific goals or set aside cash at home.'

$ws.Cells.Item(14, 1).Value = 'Citadele'
$ws.Cells.Item(14, 2).Value = 'Economic Growth and Development in the Baltic States'
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 'Citadele lowered GDP growth forecasts for all three Baltic states due to global tariff risks.'

$ws.Cells.Item(15, 1).Value = 'Citadele'
$ws.Cells.Item(15, 2).Value = 'Economic Growth and Development in the Baltic States'
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).Value = 'Citadele revised Baltic growth projections downward, expecting slower GDP expansion in Lithuania and Latvia this year.'

$ws.Cells.Item(16, 1).Value = 'Luminor'
$ws.Cells.Item(16, 2).Value = 'Business Development and Economic Challenges'
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 'Luminor ranked highest in Lithuania’s top 500 business leaders list by Verslo žinios.'

$ws.Cells.Item(17, 1).Value = 'Luminor'
$ws.Cells.Item(17, 2).Value = 'Business Development and Economic Challenges'
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 'From July 1, Luminor will change some service conditions and fees, including higher ATM withdrawal charges at Perlas terminals.'

$ws.Cells.Item(18, 1).Value = 'Luminor'
$ws.Cells.Item(18, 2).Value = 'Business Development and Economic Challenges'
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 'Luminor issued a €1.65m loan to Vilnius Public Transport for four electric passenger boats.'

$ws.Cells.Item(19, 1).Value = 'Luminor'
$ws.Cells.Item(19, 2).Value = 'Fraud and Security Issues'
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 'Luminor’s online banking and mobile app suffered outages on Tuesday evening, leaving clients unable to access services.'

$ws.Cells.Item(20, 1).Value = 'Luminor'
$ws.Cells.Item(20, 2).Value = 'Fraud and Security Issues'
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 'A second report confirmed Luminor’s e-banking and mobile services were down, with the bank apologizing for disruptions.'

$ws.Cells.Item(21, 1).Value = 'Luminor'
$ws.Cells.Item(21, 2).Value = 'Fraud and Security Issues'
$ws.Cells.Item(21, 3).Value = 3
$ws.Cells.Item(21, 4).Value = 'An attempted explosion targeted a Luminor ATM in Kelmė; the device was damaged but no injuries occurred.'

$ws.Cells.Item(22, 1).Value = 'Luminor'
$ws.Cells.Item(22, 2).Value = 'Investment Strategies and Financial Management'
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 'Luminor’s Baltic pension funds invested €18.35m into an Invalda INVL investment fund.'

$ws.Cells.Item(23, 1).Value = 'Luminor'
$ws.Cells.Item(23, 2).Value = 'Investment Strategies and Financial Management'
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 'Luminor and Revolut announced fee cuts for retail investors, while other banks held off similar moves.'

$ws.Cells.Item(24, 1).Value = 'Luminor'
$ws.Cells.Item(24, 2).Value = 'Investment Strategies and Financial Management'
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 4).Value = 'Luminor issued €300m in four-year preferred bonds with a 3.551% annual coupon, strengthening its capital base.'

$ws.Cells.Item(25, 1).Value = 'SEB'
$ws.Cells.Item(25, 2).Value = 'Banking Sector Developments and Financial Services'
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 'SEB Lithuania accepted a €490k fine from the ECB and will not appeal the decision.'

$ws.Cells.Item(26, 1).Value = 'SEB'
$ws.Cells.Item(26, 2).Value = 'Banking Sector Developments and Financial Services'
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 4).Value = 'SEB reported €146.5m net profit for H1 2025, down 6% year-on-year, with assets up 11% to €14.9bn.'

$ws.Cells.Item(27, 1).Value = 'SEB'
$ws.Cells.Item(27, 2).Value = 'Economic Growth and Forecasts'
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = 'SEB predicts GDP growth of 3% in the next half year, with a risk of a slower second half a year due to US tariffs.'

$ws.Cells.Item(28, 1).Value = 'SEB'
$ws.Cells.Item(28, 2).Value = 'Economic Growth and Forecasts'
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 4).Value = 'SEB cut Lithuania’s 2025 GDP growth forecast from 2.8% to 2.7%, projecting 2.5% growth in 2026.'

$ws.Cells.Item(29, 1).Value = 'SEB'
$ws.Cells.Item(29, 2).Value = 'Organizational Changes and Leadership in Financial Institutions'
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 'SEB''s expansion to Latvia, opening a branch there which will have positive effects for its customers.'

$ws.Cells.Item(30, 1).Value = 'SEB'
$ws.Cells.Item(30, 2).Value = 'Organizational Changes and Leadership in Financial Institutions'
$ws.Cells.Item(30, 3).Value = 2
$ws.Cells.Item(30, 4).Value = 'SEB announed to merge the 3 Baltic branches with a headquarter in Estonia and branches in Lithuania and Latvia.'

$ws.Cells.Item(31, 1).Value = 'Swedbank'
$ws.Cells.Item(31, 2).Value = 'Fraud and Security Issues'
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 'A man threw a brick at a Swedbank branch window in Vilnius; police are searching for the suspect.'

$ws.Cells.Item(32, 1).Value = 'Swedbank'
$ws.Cells.Item(32, 2).Value = 'Fraud and Security Issues'
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = 'Fraudsters extracted about €400k from Lithuanian residents in one week, with victims including a Swedbank client in Panevėžys.'

$ws.Cells.Item(33, 1).Value = 'Swedbank'
$ws.Cells.Item(33, 2).Value = 'Fraud and Security Issues'
$ws.Cells.Item(33, 3).Value = 3
$ws.Cells.Item(33, 4).Value = 'Another case reported: scammers defrauded clients of large sums, including through Swedbank accounts.'

$ws.Cells.Item(34, 1).Value = 'Swedbank'
$ws.Cells.Item(34, 2).Value = 'Customer Experience and Service Issues'
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 'A Swedbank client reported online banking outages, preventing money transfers for hours.'

$ws.Cells.Item(35, 1).Value = 'Swedbank'
$ws.Cells.Item(35, 2).Value = 'Customer Experience and Service Issues'
$ws.Cells.Item(35, 3).Value = 2
$ws.Cells.Item(35, 4).Value = 'Swedbank stated that while most services were restored after Friday’s outage, some customers still faced disruptions.'

$ws.Cells.Item(36, 1).Value = 'Swedbank'
$ws.Cells.Item(36, 2).Value = 'Economic Growth and Trends'
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = 'Swedbank reported that in Q1, Vilnius residents could afford an average apartment of 62.6 sq. m., up 11% from last quarter.'

$ws.Cells.Item(37, 1).Value = 'Swedbank'
$ws.Cells.Item(37, 2).Value = 'Economic Growth and Trends'
$ws.Cells.Item(37, 3).Value = 2
$ws.Cells.Item(37, 4).Value = 'Swedbank said the number of its retail investors in Lithuania grew nearly 70% in one year, reaching almost 100,000.'

# Update the saved selection/scroll position to match the author's final cursor location
$null = $ws.Range("D12").Select()
